$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.089.31'
$ws.Range('E2').Value = '  -1.74%  '
$ws.Range('D3').Value = '1.790.61'
$ws.Range('E3').Value = '  -0.15%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = "'222.25"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.69%  '
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = "'32.58"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.48%  '
$ws.Range('E9').Value = '  +0.27%  '
$ws.Range('D11').Value = "'0.0929"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.62%  '
$ws.Range('D12').Value = '2.046.25'
$ws.Range('E12').Value = '  -0.24%  '
$ws.Range('D13').Value = '1.794.54'
$ws.Range('E13').Value = '  -0.24%  '
$ws.Range('D14').Value = "'10.92"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.65%  '
$ws.Range('D15').Value = "'0.627"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.09%  '
$ws.Range('D16').Value = '34.061.91'
$ws.Range('E16').Value = '  -1.89%  '
$ws.Range('D17').Value = "'4.17"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -3.04%  '
$ws.Range('D18').Value = "'68.01"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.62%  '
$ws.Range('D19').Value = "'244.42"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -3.83%  '
$ws.Range('E20').Value = '  -3.19%  '
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('D22').Value = "'10.80"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.53%  '
$ws.Range('D23').Value = "'4.09"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.68%  '
$ws.Range('D24').Value = "'2.11"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.18%  '
$ws.Range('D25').Value = "'158.12"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.59%  '
$ws.Range('D26').Value = "'16.38"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('D27').Value = "'7.07"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.65%  '
$ws.Range('D28').Value = "'0.113"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.45%  '
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('D30').Value = "'0.0521"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.60%  '
$ws.Range('E31').Value = '  +0.67%  '
$ws.Range('E32').Value = '  -3.17%  '
$ws.Range('D33').Value = "'3.50"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -3.14%  '
$ws.Range('E34').Value = '  -2.81%  '
$ws.Range('D35').Value = '1.396.42'
$ws.Range('E35').Value = '  -2.98%  '
$ws.Range('D36').Value = "'0.640"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.92%  '
$ws.Range('E37').Value = '  -0.35%  '
$ws.Range('E38').Value = '  -3.10%  '
$ws.Range('D39').Value = "'79.71"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -5.93%  '
$ws.Range('D40').Value = "'0.924"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.28%  '
$ws.Range('E41').Value = '  +0.91%  '
$ws.Range('D42').Value = "'2.72"
$ws.Range('D42').ClearFormats()
$ws.Range('E43').Value = '  +1.64%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').Value = "'107.44"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +1.64%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Value = "'5.89"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.36%  '
$ws.Range('D46').Value = "'0.0494"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.49%  '
$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').Value = "'1.05"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.06%  '
$ws.Range('D48').Value = '1.946.27'
$ws.Range('E48').Value = '  +0.01%  '
$ws.Range('D49').Value = "'12.04"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.62%  '
$ws.Range('D50').Value = "'0.999"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.06%  '
$ws.Range('E51').Value = '  +2.42%  '
